$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Rename the TOC bookmark on the chapter title from _Toc440031418 to
#    _Toc442187625 (the visible bookmarkStart id="0" / _Toc428457345 stays
#    untouched).
# ---------------------------------------------------------------------------
$titleRange = $d.Content
$titleRange.Find.Execute("LIST OF ABBREVIATIONS", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$titlePara = $titleRange.Paragraphs(1)

$titleXml = $titlePara.Range.WordOpenXML
if ($titleXml -match '(?s)(<w:p\b.*?</w:p>)') {
    $titlePXml = $Matches[1]
} else {
    throw "could not locate chapter-title paragraph XML"
}

# The WordOpenXML getter stamps w14:paraId/w14:textId that are not present in
# the saved package; strip them back out so we don't introduce them.
$titlePXml = $titlePXml -replace ' w14:paraId="[0-9A-Fa-f]+"', ''
$titlePXml = $titlePXml -replace ' w14:textId="[0-9A-Fa-f]+"', ''
$titlePXml = $titlePXml -replace '_Toc440031418', '_Toc442187625'

$titlePackage = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $titlePXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$titlePara.Range.InsertXML($titlePackage) | Out-Null

# ---------------------------------------------------------------------------
# 2. Add a new "CATS / Cloud Aerosol Transport System" abbreviation entry
#    right after the existing CALIPSO entry (keeps alphabetical order).
# ---------------------------------------------------------------------------
$calipsoRange = $d.Content
$calipsoRange.Find.Execute("CALIPSO", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$calipsoPara = $calipsoRange.Paragraphs(1)

$calipsoPara.Range.InsertParagraphAfter()

$afterRange = $d.Content
$afterRange.Find.Execute("CALIPSO", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$newPara = $afterRange.Paragraphs(1).Next()

$catsPXml = '<w:p><w:pPr><w:ind w:left="1701" w:hanging="1701"/><w:rPr><w:lang w:val="en-CA" w:eastAsia="en-CA"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-CA" w:eastAsia="en-CA"/></w:rPr><w:t>CATS</w:t></w:r><w:r><w:rPr><w:lang w:val="en-CA" w:eastAsia="en-CA"/></w:rPr><w:tab/></w:r><w:r><w:t>Cloud Aerosol Transport System</w:t></w:r></w:p>'
$catsPackage = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $catsPXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$newPara.Range.InsertXML($catsPackage) | Out-Null
